$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the PARTICULARS entry on row 270 (was "SL(2-0-0)")
$ws.Range("B270").ClearContents()

# Record 3 days of Absence Undertime W/ Pay on row 271
$ws.Range("D271").Value = 3

# Remove the hidden helper/snapshot row 454; everything below shifts up by one.
$ws.Rows(454).Delete()

# Table1 regains its AutoFilter dropdown row on save.
$lo = $ws.ListObjects.Item("Table1")
$lo.ShowAutoFilter = $true
